$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values stay as literal text (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.540.10'
$ws.Range("E2").Value = '  +2.15%  '

$ws.Range("D3").Value = '2.379.58'
$ws.Range("E3").Value = '  +6.69%  '

$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").Value = '318.72'
$ws.Range("E5").Value = '  +7.50%  '

$ws.Range("D6").Value = '108.04'
$ws.Range("E6").Value = '  -3.31%  '

$ws.Range("E7").Value = '  +2.03%  '

$ws.Range("E8").Value = '  -0.20%  '

$ws.Range("D9").Value = '0.637'
$ws.Range("E9").Value = '  +4.20%  '

$ws.Range("D10").Value = '42.54'
$ws.Range("E10").Value = '  -4.67%  '

$ws.Range("D11").Value = '0.0940'
$ws.Range("E11").Value = '  +2.22%  '

$ws.Range("D12").Value = '8.72'
$ws.Range("E12").Value = '  -1.51%  '

$ws.Range("E13").Value = '  +2.30%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '0.106'
$ws.Range("E14").Value = '  +2.06%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '16.60'
$ws.Range("E15").Value = '  +9.79%  '

$ws.Range("D16").Value = '2.739.43'
$ws.Range("E16").Value = '  +6.89%  '

$ws.Range("D17").Value = '2.438.15'
$ws.Range("E17").Value = '  +8.51%  '

$ws.Range("D18").Value = '43.558.14'
$ws.Range("E18").Value = '  +2.33%  '

$ws.Range("E19").Value = '  +2.85%  '

$ws.Range("D20").Value = '7.27'
$ws.Range("E20").Value = '  -1.32%  '

$ws.Range("D21").Value = '75.33'
$ws.Range("E21").Value = '  +3.23%  '

$ws.Range("D22").Value = '3.49'
$ws.Range("E22").Value = '  -0.38%  '

$ws.Range("E23").Value = '  +5.32%  '

$ws.Range("D24").Value = '261.54'

$ws.Range("D25").Value = '9.24'
$ws.Range("E25").Value = '  +0.21%  '

$ws.Range("D26").Value = '12.04'

$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").Value = '38.90'
$ws.Range("E28").Value = '  +0.57%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '22.92'
$ws.Range("E29").Value = '  +8.65%  '

$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("D31").Value = '3.21'
$ws.Range("E31").Value = '  -1.44%  '

$ws.Range("D32").Value = '174.18'
$ws.Range("E32").Value = '  +0.32%  '

$ws.Range("E33").Value = '  +2.84%  '

$ws.Range("E34").Value = '  +3.71%  '

$ws.Range("E35").Value = '  +4.00%  '

$ws.Range("E36").Value = '  -4.02%  '

$ws.Range("E37").Value = '  -1.07%  '

$ws.Range("D38").Value = '4.08'
$ws.Range("E38").Value = '  -5.95%  '

$ws.Range("D39").Value = '0.106'
$ws.Range("E39").Value = '  +1.60%  '

$ws.Range("D40").Value = '2.83'
$ws.Range("E40").Value = '  +17.15%  '

$ws.Range("E41").Value = '  +13.55%  '

$ws.Range("D42").Value = '71.92'
$ws.Range("E42").Value = '  -0.93%  '

$ws.Range("D43").Value = '0.232'
$ws.Range("E43").Value = '  -1.54%  '

$ws.Range("E44").Value = '  +0.08%  '

$ws.Range("D45").Value = '12.60'
$ws.Range("E45").Value = '  -1.47%  '

$ws.Range("D46").Value = '5.63'
$ws.Range("E46").Value = '  +3.23%  '

$ws.Range("D47").Value = '112.68'
$ws.Range("E47").Value = '  +8.96%  '

$ws.Range("D48").Value = '9.34'
$ws.Range("E48").Value = '  +8.71%  '

$ws.Range("E49").Value = '  -0.89%  '

$ws.Range("E50").Value = '  +2.81%  '

$ws.Range("E51").Value = '  +7.85%  '
